$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column A entirely - it only held row-index numbers with header
# styling but no header text; everything currently in B:F shifts left to A:E.
$ws.Range("A1").EntireColumn.Delete()

# Fix the "MODEL_CONDITION" header text (shared string) -> "MODELCONDITION"
# After the column shift, this header now lives in column D.
$ws.Range("D1").Value = "MODELCONDITION"
